$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("21TRD09246", "Bunner", "1ST SPEED IN 1 YR >70MPH", "4511.21D4", "No Data", "Guilty", "Guilty", "`$ 0", "`$ 0", "None", "None"),
    @("21TRD09246", "Bunner", "Driving Under Suspension FTA, Fines or Child Support", "4510.111 ", "Unclassified Misdemeanor", "Guilty", "Guilty", "`$ 0", "`$ 0", "None", "None"),
    @("21TRD09246", "Bunner", "Failure to Move Over/Exercise Due Care for Public Safety, Emergency or Road Service Vehicle", "4511.213", "Minor Misdemeanor", "Guilty", "Guilty", "`$ 0", "`$ 0", "None", "None"),
    @("21TRD09246", "Bunner", "Driving Under Suspension FTA, Fines or Child Support", "4510.111 ", "Unclassified Misdemeanor", "Guilty", "Guilty", "`$ 0", "`$ 0", "None", "None"),
    @("21CRB01291", "Hemmeter", "PERMISSION REQ'D TO USE LICENSED DOCK - AMENDED to Driving Under Suspension", "1501:46-12-04", "MM", "Guilty", "Guilty", "`$ 50", "`$ 25", "10", "5"),
    @("21CRB01291", "Hemmeter", "PERMISSION REQ'D TO USE LICENSED DOCK - AMENDED to Driving Under Suspension", "1501:46-12-04", "MM", "No Contest", "Guilty", "`$ 50", "`$ 25", "10", "5"),
    @("21CRB01291", "Hemmeter", "No Operator License - Never Held", "4510.12(C)(1)", "Unclassified Misdemeanor", "No Contest", "Not Guilty", "`$ 0", "`$ 0", "None", "None"),
    @("21CRB01291", "Hemmeter", "Failure to Reinstate", "4510.21", "Unclassified Misdemeanor", "No Contest", "Guilty", "`$ 0", "`$ 0", "None", "None"),
)

$startRow = 561
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    $rng = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 11))
    $rng.NumberFormat = "@"
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowVals[$c]
    }
}

Write-Output "Added rows 561-568"